# Updates cryptos list values (Price and Volume(1h) columns) to refreshed
# scrape data, including a row-order swap for EnergySwap/PaxDollar (rows 46-47).
#
# Values are written via a throwaway ="..." formula + Copy/PasteSpecial(xlPasteValues)
# round-trip instead of a direct .Value assignment: several Price-column strings
# (e.g. "307.74", "10.20", "0.06180") are numeric-looking, and a direct .Value
# assignment lets Excel coerce them to real numbers (dropping trailing zeros,
# applying float rounding) and/or stamp the cell with a new number-format style.
# Routing them through a text formula first guarantees they land back as plain
# text with the original formatting intact and no style churn.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '="27.282.52"'
$ws.Range('D2').Copy()
$ws.Range('D2').PasteSpecial(-4163)
$ws.Range('E2').Formula = '="  +0.31%  "'
$ws.Range('E2').Copy()
$ws.Range('E2').PasteSpecial(-4163)
$ws.Range('D3').Formula = '="1.910.33"'
$ws.Range('D3').Copy()
$ws.Range('D3').PasteSpecial(-4163)
$ws.Range('E3').Formula = '="  +0.40%  "'
$ws.Range('E3').Copy()
$ws.Range('E3').PasteSpecial(-4163)
$ws.Range('E4').Formula = '="  +0.09%  "'
$ws.Range('E4').Copy()
$ws.Range('E4').PasteSpecial(-4163)
$ws.Range('D5').Formula = '="307.74"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Formula = '="  -0.14%  "'
$ws.Range('E5').Copy()
$ws.Range('E5').PasteSpecial(-4163)
$ws.Range('D7').Formula = '="0.5315"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Formula = '="  +2.02%  "'
$ws.Range('E7').Copy()
$ws.Range('E7').PasteSpecial(-4163)
$ws.Range('D8').Formula = '="0.3824"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Formula = '="  +1.58%  "'
$ws.Range('E8').Copy()
$ws.Range('E8').PasteSpecial(-4163)
$ws.Range('D9').Formula = '="0.07303"'
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Formula = '="  +0.36%  "'
$ws.Range('E9').Copy()
$ws.Range('E9').PasteSpecial(-4163)
$ws.Range('D10').Formula = '="22.12"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Formula = '="  +4.65%  "'
$ws.Range('E10').Copy()
$ws.Range('E10').PasteSpecial(-4163)
$ws.Range('D11').Formula = '="0.9025"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Formula = '="  -0.31%  "'
$ws.Range('E11').Copy()
$ws.Range('E11').PasteSpecial(-4163)
$ws.Range('D12').Formula = '="0.08191"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Formula = '="  -0.81%  "'
$ws.Range('E12').Copy()
$ws.Range('E12').PasteSpecial(-4163)
$ws.Range('D13').Formula = '="96.02"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Formula = '="  -0.88%  "'
$ws.Range('E13').Copy()
$ws.Range('E13').PasteSpecial(-4163)
$ws.Range('D14').Formula = '="5.358"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Formula = '="  +1.24%  "'
$ws.Range('E14').Copy()
$ws.Range('E14').PasteSpecial(-4163)
$ws.Range('E15').Formula = '="  -0.01%  "'
$ws.Range('E15').Copy()
$ws.Range('E15').PasteSpecial(-4163)
$ws.Range('D16').Formula = '="0.000008658"'
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Formula = '="  -0.22%  "'
$ws.Range('E16').Copy()
$ws.Range('E16').PasteSpecial(-4163)
$ws.Range('D17').Formula = '="14.81"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Formula = '="  +1.65%  "'
$ws.Range('E17').Copy()
$ws.Range('E17').PasteSpecial(-4163)
$ws.Range('E18').Formula = '="  +0.12%  "'
$ws.Range('E18').Copy()
$ws.Range('E18').PasteSpecial(-4163)
$ws.Range('D19').Formula = '="1.283.99"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Formula = '="  -32.43%  "'
$ws.Range('E19').Copy()
$ws.Range('E19').PasteSpecial(-4163)
$ws.Range('D20').Formula = '="27.317.93"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Formula = '="  +0.30%  "'
$ws.Range('E20').Copy()
$ws.Range('E20').PasteSpecial(-4163)
$ws.Range('D21').Formula = '="5.071"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Formula = '="  -0.52%  "'
$ws.Range('E21').Copy()
$ws.Range('E21').PasteSpecial(-4163)
$ws.Range('E22').Formula = '="  +1.54%  "'
$ws.Range('E22').Copy()
$ws.Range('E22').PasteSpecial(-4163)
$ws.Range('D23').Formula = '="6.521"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Formula = '="  +1.32%  "'
$ws.Range('E23').Copy()
$ws.Range('E23').PasteSpecial(-4163)
$ws.Range('D24').Formula = '="149.90"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Formula = '="  +2.35%  "'
$ws.Range('E24').Copy()
$ws.Range('E24').PasteSpecial(-4163)
$ws.Range('D25').Formula = '="2.294"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Formula = '="  -1.39%  "'
$ws.Range('E25').Copy()
$ws.Range('E25').PasteSpecial(-4163)
$ws.Range('D26').Formula = '="18.25"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Formula = '="  +0.11%  "'
$ws.Range('E26').Copy()
$ws.Range('E26').PasteSpecial(-4163)
$ws.Range('D27').Formula = '="1.742"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Formula = '="  -0.35%  "'
$ws.Range('E27').Copy()
$ws.Range('E27').PasteSpecial(-4163)
$ws.Range('D28').Formula = '="116.77"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Formula = '="  +1.52%  "'
$ws.Range('E28').Copy()
$ws.Range('E28').PasteSpecial(-4163)
$ws.Range('D29').Formula = '="4.832"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Formula = '="  -0.05%  "'
$ws.Range('E29').Copy()
$ws.Range('E29').PasteSpecial(-4163)
$ws.Range('D30').Formula = '="4.822"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Formula = '="  -1.67%  "'
$ws.Range('E30').Copy()
$ws.Range('E30').PasteSpecial(-4163)
$ws.Range('D31').Formula = '="0.09283"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Formula = '="  +0.16%  "'
$ws.Range('E31').Copy()
$ws.Range('E31').PasteSpecial(-4163)
$ws.Range('D32').Formula = '="0.8334"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Formula = '="  +4.22%  "'
$ws.Range('E32').Copy()
$ws.Range('E32').PasteSpecial(-4163)
$ws.Range('D33').Formula = '="0.05071"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Formula = '="  -0.27%  "'
$ws.Range('E33').Copy()
$ws.Range('E33').PasteSpecial(-4163)
$ws.Range('D34').Formula = '="1.229"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Formula = '="  -1.23%  "'
$ws.Range('E34').Copy()
$ws.Range('E34').PasteSpecial(-4163)
$ws.Range('D35').Formula = '="3.001"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Formula = '="  +1.85%  "'
$ws.Range('E35').Copy()
$ws.Range('E35').PasteSpecial(-4163)
$ws.Range('D36').Formula = '="3.358"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Formula = '="  -1.97%  "'
$ws.Range('E36').Copy()
$ws.Range('E36').PasteSpecial(-4163)
$ws.Range('D37').Formula = '="2.685"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Formula = '="  +3.40%  "'
$ws.Range('E37').Copy()
$ws.Range('E37').PasteSpecial(-4163)
$ws.Range('D38').Formula = '="0.5756"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Formula = '="  +0.71%  "'
$ws.Range('E38').Copy()
$ws.Range('E38').PasteSpecial(-4163)
$ws.Range('D39').Formula = '="0.02007"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Formula = '="  +0.40%  "'
$ws.Range('E39').Copy()
$ws.Range('E39').PasteSpecial(-4163)
$ws.Range('D40').Formula = '="1.077"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('D41').Formula = '="9.402"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Formula = '="  +4.23%  "'
$ws.Range('E41').Copy()
$ws.Range('E41').PasteSpecial(-4163)
$ws.Range('D42').Formula = '="6.567"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Formula = '="  -0.40%  "'
$ws.Range('E42').Copy()
$ws.Range('E42').PasteSpecial(-4163)
$ws.Range('D43').Formula = '="117.01"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Formula = '="  +0.07%  "'
$ws.Range('E43').Copy()
$ws.Range('E43').PasteSpecial(-4163)
$ws.Range('D44').Formula = '="0.1523"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Formula = '="  +0.33%  "'
$ws.Range('E44').Copy()
$ws.Range('E44').PasteSpecial(-4163)
$ws.Range('D45').Formula = '="0.4922"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Formula = '="  +1.37%  "'
$ws.Range('E45').Copy()
$ws.Range('E45').PasteSpecial(-4163)
$ws.Range('B46').Formula = '="EnergySwap"'
$ws.Range('B46').Copy()
$ws.Range('B46').PasteSpecial(-4163)
$ws.Range('C46').Formula = '="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"'
$ws.Range('C46').Copy()
$ws.Range('C46').PasteSpecial(-4163)
$ws.Range('D46').Formula = '="10.20"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Formula = '="  +1.11%  "'
$ws.Range('E46').Copy()
$ws.Range('E46').PasteSpecial(-4163)
$ws.Range('B47').Formula = '="PaxDollar"'
$ws.Range('B47').Copy()
$ws.Range('B47').PasteSpecial(-4163)
$ws.Range('C47').Formula = '="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"'
$ws.Range('C47').Copy()
$ws.Range('C47').PasteSpecial(-4163)
$ws.Range('D47').Formula = '="1.001"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Formula = '="  +0.03%  "'
$ws.Range('E47').Copy()
$ws.Range('E47').PasteSpecial(-4163)
$ws.Range('D48').Formula = '="1.639"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Formula = '="  +0.57%  "'
$ws.Range('E48').Copy()
$ws.Range('E48').PasteSpecial(-4163)
$ws.Range('D49').Formula = '="38.83"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Formula = '="  +2.93%  "'
$ws.Range('E49').Copy()
$ws.Range('E49').PasteSpecial(-4163)
$ws.Range('D50').Formula = '="0.06180"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Formula = '="  +3.76%  "'
$ws.Range('E50').Copy()
$ws.Range('E50').PasteSpecial(-4163)
$ws.Range('E51').Formula = '="  -0.34%  "'
$ws.Range('E51').Copy()
$ws.Range('E51').PasteSpecial(-4163)

$excel.CutCopyMode = 0
Write-Host "Applied 96 cell updates"
